$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("股票")

# ---------------------------------------------------------------------------
# Target layout change (sheet "股票" / stock):
#   insert a new "category" column right after "property_category" (old col I
#   shifts date/legislator_name/legislator_id one column to the right), and
#   append two new trailing columns "source_file" / "index" right after
#   "legislator_id".
#
# Because EntireColumn.Insert() is the only operation in this COM surface
# that reliably carries the existing cell style (header style 1 / data style
# 2) across a shift, we do all three inserts that way and then only move
# *values* (never styles) into their final homes.
# ---------------------------------------------------------------------------

# 1) Insert the "category" column before the existing "date" column (old I).
#    Old I/J/K (date/legislator_name/legislator_id) shift right to J/K/L.
$ws.Range("I1").EntireColumn.Insert()

# 2) Insert two more columns before "legislator_id" (now at L), so
#    legislator_id is pushed to N, leaving two freshly-styled blank columns
#    at L and M for "source_file" / "index".
$ws.Range("L1").EntireColumn.Insert()
$ws.Range("L1").EntireColumn.Insert()

# 3) Header row: move the "legislator_id" header text back to L1 (its final
#    column), then label the newly created columns.
$legIdHeader = $ws.Range("N1").Value2
$ws.Range("L1").Value = $legIdHeader
$ws.Range("I1").Value = "category"
$ws.Range("M1").Value = "source_file"
$ws.Range("N1").Value = "index"

# 4) Data rows: same value shuffle, plus the new per-row values.
for ($r = 2; $r -le 8; $r++) {
    $idxValue = $ws.Range("A$r").Value2
    $legIdValue = $ws.Range("N$r").Value2

    $ws.Range("I$r").Value = "normal"
    $ws.Range("L$r").Value = $legIdValue
    $ws.Range("M$r").Value = "tmp581f1"
    $ws.Range("N$r").Value = $idxValue
}
